$d = $word.ActiveDocument

# --- Change 1: merge the three "list" / " etc. here " / "Counter" runs
# into a single run. They already read as one contiguous phrase, so a
# Find/Replace of the phrase with itself collapses the (formatting-
# identical) adjacent runs back into a single <w:r>.
$d.Content.Find.Execute("list etc. here Counter", $true, $false, $false, `
    $false, $false, $true, 1, $false, "list etc. here Counter", 2) | Out-Null

# --- Change 2: append two new paragraphs (a "Stack ..." line and a
# trailing blank line, both in Times New Roman 12pt) right before the
# document's final empty paragraph.
$lastPara = $d.Paragraphs.Last
$insertionRange = $lastPara.Range

$xml = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" ' + `
    'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
    'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" ' + `
    'w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/>' + `
    '<w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" ' + `
    'w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' + `
    '<w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Stack ' + [char]0x2013 + `
    ' can be implemented as a linked list and as an array.</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" ' + `
    'w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/>' + `
    '<w:szCs w:val="24"/></w:rPr></w:pPr></w:p>' + `
    '<w:p w14:paraId="29F3377B" w14:textId="77777777" w:rsidR="00762A15" ' + `
    'w:rsidRDefault="00762A15" w:rsidP="00762A15"/></pkg:xmlData>'

$insertionRange.InsertXML($xml)
